$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column B to text so the dd/mm/yyyy strings are not
# auto-converted into Excel date serials when the .Value is assigned below.
$ws.Range("B175:B185").NumberFormat = "@"

# Row 175
$ws.Cells.Item(175, 1).Value = 173
$ws.Cells.Item(175, 2).Value = "12/06/2024"
$ws.Cells.Item(175, 3).Value = 147.2666666666667
$ws.Cells.Item(175, 4).Value = 149
$ws.Cells.Item(175, 5).Value = 143
$ws.Cells.Item(175, 6).Value = 139
$ws.Cells.Item(175, 7).Value = 64.75
$ws.Cells.Item(175, 8).Value = 148.4
$ws.Cells.Item(175, 9).Value = 67.83799999999999
$ws.Cells.Item(175, 10).Value = 62
$ws.Cells.Item(175, 11).Value = 133.3725
$ws.Cells.Item(175, 12).Value = 170.957
$ws.Cells.Item(175, 13).Value = 120
$ws.Cells.Item(175, 14).Value = 196.5
$ws.Cells.Item(175, 15).Value = 182
$ws.Cells.Item(175, 16).Value = 179
$ws.Cells.Item(175, 17).Value = 175.5
$ws.Cells.Item(175, 18).Value = 99
$ws.Cells.Item(175, 19).Value = 185
$ws.Cells.Item(175, 20).Value = 0.3192307692307693
$ws.Cells.Item(175, 21).Value = 64.48999999999999
$ws.Cells.Item(175, 22).Value = 118.25
$ws.Cells.Item(175, 23).Value = 64.48999999999999

# Row 176
$ws.Cells.Item(176, 1).Value = 174
$ws.Cells.Item(176, 2).Value = "13/06/2024"
$ws.Cells.Item(176, 3).Value = 147.1538461538462
$ws.Cells.Item(176, 4).Value = 146.5
$ws.Cells.Item(176, 5).Value = 143
$ws.Cells.Item(176, 6).Value = 140
$ws.Cells.Item(176, 7).Value = 64.75
$ws.Cells.Item(176, 8).Value = 150
$ws.Cells.Item(176, 9).Value = 67.83799999999999
$ws.Cells.Item(176, 10).Value = 62
$ws.Cells.Item(176, 11).Value = 131.0225
$ws.Cells.Item(176, 12).Value = 167.5376086956522
$ws.Cells.Item(176, 13).Value = 120
$ws.Cells.Item(176, 14).Value = 185
$ws.Cells.Item(176, 15).Value = 182
$ws.Cells.Item(176, 16).Value = 179
$ws.Cells.Item(176, 17).Value = 175.5
$ws.Cells.Item(176, 18).Value = 99
$ws.Cells.Item(176, 19).Value = 185
$ws.Cells.Item(176, 20).Value = 0.3192307692307693
$ws.Cells.Item(176, 21).Value = 64.48999999999999
$ws.Cells.Item(176, 22).Value = 118.25
$ws.Cells.Item(176, 23).Value = 64.48999999999999

# Row 177
$ws.Cells.Item(177, 1).Value = 175
$ws.Cells.Item(177, 2).Value = "14/06/2024"
$ws.Cells.Item(177, 3).Value = 151.2584375
$ws.Cells.Item(177, 4).Value = 150.4035714285714
$ws.Cells.Item(177, 5).Value = 144.1666666666667
$ws.Cells.Item(177, 6).Value = 141
$ws.Cells.Item(177, 7).Value = 64.75
$ws.Cells.Item(177, 8).Value = 122.75
$ws.Cells.Item(177, 9).Value = 67.83799999999999
$ws.Cells.Item(177, 10).Value = 62
$ws.Cells.Item(177, 11).Value = 138.31125
$ws.Cells.Item(177, 12).Value = 173.1102985074627
$ws.Cells.Item(177, 13).Value = 120
$ws.Cells.Item(177, 14).Value = 185
$ws.Cells.Item(177, 15).Value = 182
$ws.Cells.Item(177, 16).Value = 179
$ws.Cells.Item(177, 17).Value = 175.5
$ws.Cells.Item(177, 18).Value = 99
$ws.Cells.Item(177, 19).Value = 185
$ws.Cells.Item(177, 20).Value = 0.3192307692307693
$ws.Cells.Item(177, 21).Value = 64.48999999999999
$ws.Cells.Item(177, 22).Value = 118.25
$ws.Cells.Item(177, 23).Value = 64.48999999999999

# Row 178
$ws.Cells.Item(178, 1).Value = 176
$ws.Cells.Item(178, 2).Value = "17/06/2024"
$ws.Cells.Item(178, 3).Value = 156.4611111111111
$ws.Cells.Item(178, 4).Value = 152.2
$ws.Cells.Item(178, 5).Value = 144.5
$ws.Cells.Item(178, 6).Value = 143
$ws.Cells.Item(178, 7).Value = 64.75
$ws.Cells.Item(178, 8).Value = 122.75
$ws.Cells.Item(178, 9).Value = 67.83799999999999
$ws.Cells.Item(178, 10).Value = 62
$ws.Cells.Item(178, 11).Value = 146.0454545454545
$ws.Cells.Item(178, 12).Value = 180.1802941176471
$ws.Cells.Item(178, 13).Value = 120
$ws.Cells.Item(178, 14).Value = 185
$ws.Cells.Item(178, 15).Value = 182
$ws.Cells.Item(178, 16).Value = 179
$ws.Cells.Item(178, 17).Value = 175.5
$ws.Cells.Item(178, 18).Value = 99
$ws.Cells.Item(178, 19).Value = 196
$ws.Cells.Item(178, 20).Value = 0.3192307692307693
$ws.Cells.Item(178, 21).Value = 64.48999999999999
$ws.Cells.Item(178, 22).Value = 118.25
$ws.Cells.Item(178, 23).Value = 64.48999999999999

# Row 179
$ws.Cells.Item(179, 1).Value = 177
$ws.Cells.Item(179, 2).Value = "18/06/2024"
$ws.Cells.Item(179, 3).Value = 157.6777777777778
$ws.Cells.Item(179, 4).Value = 152.5833333333333
$ws.Cells.Item(179, 5).Value = 146
$ws.Cells.Item(179, 6).Value = 143
$ws.Cells.Item(179, 7).Value = 64.75
$ws.Cells.Item(179, 8).Value = 164
$ws.Cells.Item(179, 9).Value = 67.83799999999999
$ws.Cells.Item(179, 10).Value = 62
$ws.Cells.Item(179, 11).Value = 146.0454545454545
$ws.Cells.Item(179, 12).Value = 181.1875
$ws.Cells.Item(179, 13).Value = 120
$ws.Cells.Item(179, 14).Value = 185
$ws.Cells.Item(179, 15).Value = 182
$ws.Cells.Item(179, 16).Value = 179
$ws.Cells.Item(179, 17).Value = 175.5
$ws.Cells.Item(179, 18).Value = 99
$ws.Cells.Item(179, 19).Value = 200
$ws.Cells.Item(179, 20).Value = 0.3192307692307693
$ws.Cells.Item(179, 21).Value = 64.48999999999999
$ws.Cells.Item(179, 22).Value = 118.25
$ws.Cells.Item(179, 23).Value = 64.48999999999999

# Row 180
$ws.Cells.Item(180, 1).Value = 178
$ws.Cells.Item(180, 2).Value = "19/06/2024"
$ws.Cells.Item(180, 3).Value = 155.6666666666667
$ws.Cells.Item(180, 4).Value = 152.5833333333333
$ws.Cells.Item(180, 5).Value = 146
$ws.Cells.Item(180, 6).Value = 143
$ws.Cells.Item(180, 7).Value = 64.75
$ws.Cells.Item(180, 8).Value = 164
$ws.Cells.Item(180, 9).Value = 67.83799999999999
$ws.Cells.Item(180, 10).Value = 62
$ws.Cells.Item(180, 11).Value = 146.0454545454545
$ws.Cells.Item(180, 12).Value = 177.555
$ws.Cells.Item(180, 13).Value = 120
$ws.Cells.Item(180, 14).Value = 185
$ws.Cells.Item(180, 15).Value = 182
$ws.Cells.Item(180, 16).Value = 179
$ws.Cells.Item(180, 17).Value = 175.5
$ws.Cells.Item(180, 18).Value = 99
$ws.Cells.Item(180, 19).Value = 200
$ws.Cells.Item(180, 20).Value = 0.3192307692307693
$ws.Cells.Item(180, 21).Value = 64.48999999999999
$ws.Cells.Item(180, 22).Value = 118.25
$ws.Cells.Item(180, 23).Value = 64.48999999999999

# Row 181
$ws.Cells.Item(181, 1).Value = 179
$ws.Cells.Item(181, 2).Value = "20/06/2024"
$ws.Cells.Item(181, 3).Value = 157.2
$ws.Cells.Item(181, 4).Value = 153
$ws.Cells.Item(181, 5).Value = 146
$ws.Cells.Item(181, 6).Value = 143
$ws.Cells.Item(181, 7).Value = 64.75
$ws.Cells.Item(181, 8).Value = 161
$ws.Cells.Item(181, 9).Value = 67.83799999999999
$ws.Cells.Item(181, 10).Value = 62
$ws.Cells.Item(181, 11).Value = 145.8333333333333
$ws.Cells.Item(181, 12).Value = 179.6538461538462
$ws.Cells.Item(181, 13).Value = 120
$ws.Cells.Item(181, 14).Value = 185
$ws.Cells.Item(181, 15).Value = 182
$ws.Cells.Item(181, 16).Value = 179
$ws.Cells.Item(181, 17).Value = 175.5
$ws.Cells.Item(181, 18).Value = 99
$ws.Cells.Item(181, 19).Value = 200
$ws.Cells.Item(181, 20).Value = 0.3192307692307693
$ws.Cells.Item(181, 21).Value = 64.48999999999999
$ws.Cells.Item(181, 22).Value = 118.25
$ws.Cells.Item(181, 23).Value = 64.48999999999999

# Row 182
$ws.Cells.Item(182, 1).Value = 180
$ws.Cells.Item(182, 2).Value = "21/06/2024"
$ws.Cells.Item(182, 3).Value = 160.7992857142857
$ws.Cells.Item(182, 4).Value = 158.024375
$ws.Cells.Item(182, 5).Value = 148.5555555555555
$ws.Cells.Item(182, 6).Value = 143
$ws.Cells.Item(182, 7).Value = 64.75
$ws.Cells.Item(182, 8).Value = 168.575
$ws.Cells.Item(182, 9).Value = 67.83799999999999
$ws.Cells.Item(182, 10).Value = 62
$ws.Cells.Item(182, 11).Value = 155.5
$ws.Cells.Item(182, 12).Value = 187.3439285714286
$ws.Cells.Item(182, 13).Value = 120
$ws.Cells.Item(182, 14).Value = 196.8333333333333
$ws.Cells.Item(182, 15).Value = 182
$ws.Cells.Item(182, 16).Value = 179
$ws.Cells.Item(182, 17).Value = 175.5
$ws.Cells.Item(182, 18).Value = 99
$ws.Cells.Item(182, 19).Value = 205
$ws.Cells.Item(182, 20).Value = 0.3192307692307693
$ws.Cells.Item(182, 21).Value = 64.48999999999999
$ws.Cells.Item(182, 22).Value = 118.25
$ws.Cells.Item(182, 23).Value = 64.48999999999999

# Row 183
$ws.Cells.Item(183, 1).Value = 181
$ws.Cells.Item(183, 2).Value = "24/06/2024"
$ws.Cells.Item(183, 3).Value = 160.3684210526316
$ws.Cells.Item(183, 4).Value = 155.6666666666667
$ws.Cells.Item(183, 5).Value = 148.5555555555555
$ws.Cells.Item(183, 6).Value = 144
$ws.Cells.Item(183, 7).Value = 64.75
$ws.Cells.Item(183, 8).Value = 168.575
$ws.Cells.Item(183, 9).Value = 67.83799999999999
$ws.Cells.Item(183, 10).Value = 62
$ws.Cells.Item(183, 11).Value = 156.9166666666667
$ws.Cells.Item(183, 12).Value = 188.6785714285714
$ws.Cells.Item(183, 13).Value = 120
$ws.Cells.Item(183, 14).Value = 196.8333333333333
$ws.Cells.Item(183, 15).Value = 182
$ws.Cells.Item(183, 16).Value = 179
$ws.Cells.Item(183, 17).Value = 175.5
$ws.Cells.Item(183, 18).Value = 99
$ws.Cells.Item(183, 19).Value = 205
$ws.Cells.Item(183, 20).Value = 0.3192307692307693
$ws.Cells.Item(183, 21).Value = 64.48999999999999
$ws.Cells.Item(183, 22).Value = 118.25
$ws.Cells.Item(183, 23).Value = 64.48999999999999

# Row 184
$ws.Cells.Item(184, 1).Value = 182
$ws.Cells.Item(184, 2).Value = "25/06/2024"
$ws.Cells.Item(184, 3).Value = 161.9607142857143
$ws.Cells.Item(184, 4).Value = 156.8285714285714
$ws.Cells.Item(184, 5).Value = 147.1666666666667
$ws.Cells.Item(184, 6).Value = 143
$ws.Cells.Item(184, 7).Value = 64.75
$ws.Cells.Item(184, 8).Value = 168.575
$ws.Cells.Item(184, 9).Value = 67.83799999999999
$ws.Cells.Item(184, 10).Value = 62
$ws.Cells.Item(184, 11).Value = 155.3557142857143
$ws.Cells.Item(184, 12).Value = 190.9038333333333
$ws.Cells.Item(184, 13).Value = 120
$ws.Cells.Item(184, 14).Value = 196.8333333333333
$ws.Cells.Item(184, 15).Value = 182
$ws.Cells.Item(184, 16).Value = 179
$ws.Cells.Item(184, 17).Value = 175.5
$ws.Cells.Item(184, 18).Value = 99
$ws.Cells.Item(184, 19).Value = 210
$ws.Cells.Item(184, 20).Value = 0.3192307692307693
$ws.Cells.Item(184, 21).Value = 64.48999999999999
$ws.Cells.Item(184, 22).Value = 118.25
$ws.Cells.Item(184, 23).Value = 64.48999999999999

# Row 185
$ws.Cells.Item(185, 1).Value = 183
$ws.Cells.Item(185, 2).Value = "26/06/2024"
$ws.Cells.Item(185, 3).Value = 165.3409090909091
$ws.Cells.Item(185, 4).Value = 157.9666666666667
$ws.Cells.Item(185, 5).Value = 149.5
$ws.Cells.Item(185, 6).Value = 142.75
$ws.Cells.Item(185, 7).Value = 64.75
$ws.Cells.Item(185, 8).Value = 174.35
$ws.Cells.Item(185, 9).Value = 67.83799999999999
$ws.Cells.Item(185, 10).Value = 62
$ws.Cells.Item(185, 11).Value = 158.15
$ws.Cells.Item(185, 12).Value = 193.7390322580645
$ws.Cells.Item(185, 13).Value = 120
$ws.Cells.Item(185, 14).Value = 199.5
$ws.Cells.Item(185, 15).Value = 182
$ws.Cells.Item(185, 16).Value = 179
$ws.Cells.Item(185, 17).Value = 175.5
$ws.Cells.Item(185, 18).Value = 99
$ws.Cells.Item(185, 19).Value = 210
$ws.Cells.Item(185, 20).Value = 0.3192307692307693
$ws.Cells.Item(185, 21).Value = 64.48999999999999
$ws.Cells.Item(185, 22).Value = 118.25
$ws.Cells.Item(185, 23).Value = 64.48999999999999

# The source workbook stores these date-strings with default (General) cell
# formatting, so drop the temporary text format now that the literal strings
# are safely stored - this restores the default style (no explicit style index).
$ws.Range("B175:B185").ClearFormats()

# Copy formatting (bold, border, center alignment) from the last pre-existing
# row's A cell into the new A-column cells so they match the existing header-style column
$ws.Cells.Item(174, 1).Copy() | Out-Null
$ws.Range("A175:A185").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
